# Update column G ("K") values for rows 2-30 per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4, 2, 4, 7, 2, 5, 6, 2, 2, 10, 4, 7, 6, 3, 6, 4, 6, 6, 2, 3, 6, 4, 3, 5, 13, 4, 3, 4, 2)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
